$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.280.79"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.18"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.38"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.875.16"
$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.605.91"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.551"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.56"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.264.93"
$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.38"
$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.78"
$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.10"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("E34").Value = "  +1.71%  "

$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.139.90"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.557"
$ws.Range("E38").Value = "  +1.85%  "

$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  +2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.33"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("E44").Value = "  -1.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.783.73"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.40"
$ws.Range("E46").Value = "  +1.54%  "

$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.75"
$ws.Range("E49").Value = "  +3.08%  "

$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0975"
$ws.Range("E51").Value = "  +1.90%  "
